$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"0.992651046408802"
$ws.Range("E2").Value = [double]"0.992651046408802"

$ws.Range("D3").Value = [double]"0.99886708929831"
$ws.Range("E3").Value = [double]"0.99886708929831"

$ws.Range("D4").Value = [double]"0.9999999999265423"
$ws.Range("E4").Value = [double]"0.9999999999265423"

$ws.Range("D5").Value = [double]"5.503754144703708E-09"
$ws.Range("E5").Value = [double]"5.503754144703708E-09"

$ws.Range("D6").Value = [double]"9.0827815975469E-19"
$ws.Range("E6").Value = [double]"9.0827815975469E-19"

$ws.Range("D7").Value = [double]"6.001768742932272E-15"
$ws.Range("E7").Value = [double]"0.999999999999994"

$ws.Range("D8").Value = [double]"0.9999999543435237"
$ws.Range("E8").Value = [double]"4.56564762663092E-08"

$ws.Range("D9").Value = [double]"0.9850699954689528"
$ws.Range("E9").Value = [double]"0.0149300045310472"

$ws.Range("D10").Value = [double]"2.98732360898077E-05"
$ws.Range("E10").Value = [double]"0.9999701267639102"

$ws.Range("D11").Value = [double]"0.9999954721950703"
$ws.Range("E11").Value = [double]"4.527804929654877E-06"
$ws.Range("F11").Value = [double]"7.821080207824707"
